# Commit: "Added Graphs + edited data"
# Data correction on Sheet1: two Altitude readings (column F) are updated
# to 30000, and the sheet's active selection/scroll position is moved
# from the bottom of the data (H11, scrolled to row 5) back up to F6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Edited data: Altitude column (F) ---
$ws.Range("F2").Value = 30000
$ws.Range("F4").Value = 30000

# --- View: clear the scrolled-down top-left cell and move the selection ---
$ws.Activate() | Out-Null
$ws.Range("F6").Select() | Out-Null
